$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "315.44"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "2.54%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "39.51"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "-1.34%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.127"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "0.02%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.08170"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "0.89%"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.971"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "1.62%"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "3.83%"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "8.339"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "2.68%"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.9370"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "0.75%"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.1297"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "-7.90%"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "2.77%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.09018"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "-0.62%"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.03501"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-0.54%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.09748"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "-0.60%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.001416"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "1.55%"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.006099"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "3.62%"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "-6.15%"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.275"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-3.49%"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "1.94%"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "-0.18%"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.965"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "6.38%"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.2495"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "-0.20%"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04367"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "-0.52%"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001242"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "2.15%"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "9.62%"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0003896"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "198.58%"
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "-7.58%"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "8.54%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.05161"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "2.32%"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.007762"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "4.91%"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.01042"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "6.02%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.1404"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "2.99%"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-4.40%"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.009291"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "1.02%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00006939"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "9.01%"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.00000000752"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "0.20%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.002887"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "0.46%"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.001694"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "30.27%"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.00002105"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "0.20%"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0002005"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "0.20%"
$c.Style = "Normal"
